$d = $word.ActiveDocument

$replacements = @(
    @{old = "64×41="; new = "90×54="},
    @{old = "60×42="; new = "23×62="},
    @{old = "12×76="; new = "83×18="},
    @{old = "18×25="; new = "48×72="},
    @{old = "84×98="; new = "66×30="},
    @{old = "16×54="; new = "92×34="},
    @{old = "59×37="; new = "75×21="},
    @{old = "39×31="; new = "70×93="},
    @{old = "16×70="; new = "31×63="},
    @{old = "63×61="; new = "93×37="},
    @{old = "75×59="; new = "96×23="},
    @{old = "46×68="; new = "97×48="},
    @{old = "44×87="; new = "88×99="},
    @{old = "61×60="; new = "70×97="},
    @{old = "77×86="; new = "17×80="},
    @{old = "95×14="; new = "84×67="},
    @{old = "35×31="; new = "26×19="},
    @{old = "84×37="; new = "31×97="},
    @{old = "46×14="; new = "29×27="},
    @{old = "71×30="; new = "83×37="},
    @{old = "99×70="; new = "60×11="},
    @{old = "48×76="; new = "47×93="},
    @{old = "61×91="; new = "42×87="},
    @{old = "68×19="; new = "13×97="},
    @{old = "21×19="; new = "68×32="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
